$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the store names between row 3 and row 4 (Ponta Negra <-> Vieiralves)
$ws.Range("A3").Value = "Bibi Cell Vieiralves"
$ws.Range("A4").Value = "Bibi Cell Ponta Negra"

# Row 2 - Bibi Cell Mundi: update G/H and total
$ws.Range("G2").Value = 6867.91
$ws.Range("H2").Value = 4071
$ws.Range("AG2").Value = 82434.75999999999

# Row 3 - now Bibi Cell Vieiralves
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 4464
$ws.Range("D3").Value = 3201
$ws.Range("E3").Value = 3626
$ws.Range("F3").Value = 4238
$ws.Range("G3").Value = 4951
$ws.Range("H3").Value = 7617.5
$ws.Range("AG3").Value = 28097.5

# Row 4 - now Bibi Cell Ponta Negra
$ws.Range("B4").Value = 1800.01
$ws.Range("C4").Value = 4670
$ws.Range("D4").Value = 1748.51
$ws.Range("E4").Value = 5592
$ws.Range("F4").Value = 3002
$ws.Range("G4").Value = 823
$ws.Range("H4").Value = 3138.5
$ws.Range("AG4").Value = 20774.02

# Row 5 - Bibi Cell Manauara
$ws.Range("G5").Value = 2892
$ws.Range("H5").Value = 4208.4
$ws.Range("AG5").Value = 19215.4

# Row 6 - total
$ws.Range("G6").Value = 15533.91
$ws.Range("H6").Value = 19035.4
$ws.Range("AG6").Value = 150521.68
